$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("button_apiKeyAction_trNthChild" / "2") is removed entirely; the
# former column B ("input_KeyName" / blank) shifts left to become column A.
$ws.Columns.Item(1).Delete()
